$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text columns D and E keep their original "Text" format so that
# numeric-looking strings (e.g. "0.9999") are preserved as text, not numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.826.04"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.19%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.887.76"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.38%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7525"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.86%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.24"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.59%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9996"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3124"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.36%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.30"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.80%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07118"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.18%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08492"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +4.96%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7602"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.64%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.358"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.67%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.853.83"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.60%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.43"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.52%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.159"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.90%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.885.64"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.03%  "

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.95%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.67"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007796"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.35%  "

$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9993"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.16%  "

$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.135.53"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.29%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.022"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.61%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1595"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.34%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.378"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.66%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.42"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.87%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.75"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.10%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.25%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.493"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.95%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.536"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.44%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.79%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.125"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.59%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.62%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.241"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.12%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7485"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.49%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.21%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.709"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.95%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01946"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.42%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.775"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.73%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4457"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.39%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.109"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.62%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.090.53"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.46%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "72.39"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.65%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8554"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.35%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9994"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.06%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.713"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.51%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.47"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.06%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.861"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.51%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.043"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.43%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.048.47"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.70%  "
